$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks first - they will be re-created after the
# new rows are inserted (so they point at the correct, shifted cells).
$ws.Hyperlinks.Delete()

# Insert two new rows above row 1; this pushes the four existing member
# rows down from 1-4 to 3-6.
$ws.Rows("1:2").Insert()

# Row 2 (new) is filled in first so its e-mail address is registered before
# row 1's in the shared-string table, matching the order the author actually
# typed the test rows in.
# Row 2 (new): john_tan / John / Tan with a newly tested e-mail address.
$ws.Range("A2").Value = "john_tan"
$ws.Range("B2").Value = "John"
$ws.Range("C2").Value = "Tan"
$ws.Range("D2").Value = "kijas59503@petloca.com"

# Row 1 (new): john_tan / John / Tan with another newly tested e-mail address.
$ws.Range("A1").Value = "john_tan"
$ws.Range("B1").Value = "John"
$ws.Range("C1").Value = "Tan"
$ws.Range("D1").Value = "autumnlpx@gmail.com"

# Re-create the hyperlinks for the four original rows (now at 3-6), followed
# by the two brand-new ones (rows 2 and 1), matching the order in which the
# workbook author ended up with them.
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:louiseairahnicole@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:20007073@myrp.edu.sg")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:nclthr@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:florian.muljono@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:kijas59503@petloca.com")
$ws.Hyperlinks.Add($ws.Range("D1"), "mailto:autumnlpx@gmail.com")

# Re-apply the Hyperlink cell style to every link cell so they keep sharing
# the workbook's existing hyperlink format instead of a freshly duplicated one.
$ws.Range("D1").Style = "Hyperlink"
$ws.Range("D2").Style = "Hyperlink"
$ws.Range("D3").Style = "Hyperlink"
$ws.Range("D4").Style = "Hyperlink"
$ws.Range("D5").Style = "Hyperlink"
$ws.Range("D6").Style = "Hyperlink"

# Put the active selection on D2, as in the edited workbook.
$ws.Range("D2").Select()
